# Update COVID-19 country data: refresh case counts, re-sort several countries that
# swapped rank order due to updated figures, bump the "datos actualizados" timestamp,
# and append a new row for Papua Nueva Guinea (first reported case).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 1: Datos actualizados a 27 de Marzo de 2020 a las 00:12
$ws.Range("A1").Value = 'Datos actualizados a 27 de Marzo de 2020 a las 00:12'

# row 4
$ws.Range("B4").Value = 83206
$ws.Range("C4").Value = 14995
$ws.Range("E4").Value = 80141

# row 15
$ws.Range("F15").Value = 96

# row 17
$ws.Range("D17").Value = 228
$ws.Range("E17").Value = 3776

# row 21
$ws.Range("B21").Value = 2996
$ws.Range("C21").Value = 320
$ws.Range("E21").Value = 2813

# row 22
$ws.Range("B22").Value = 2985
$ws.Range("C22").Value = 431
$ws.Range("E22").Value = 2902

# row 31
$ws.Range("B31").Value = 1387
$ws.Range("C31").Value = 80
$ws.Range("E31").Value = 981

# row 34
$ws.Range("B34").Value = 1201
$ws.Range("C34").Value = 138
$ws.Range("E34").Value = 1171

# row 43
$ws.Range("D43").Value = 82
$ws.Range("E43").Value = 718

# row 44
$ws.Range("B44").Value = 727
$ws.Range("C44").Value = 70
$ws.Range("E44").Value = 662
$ws.Range("G44").Value = 8
$ws.Range("H44").Value = 20

# row 119: Ruanda
$ws.Range("A119").Value = 'Ruanda'
$ws.Range("B119").Value = 50
$ws.Range("C119").Value = 9
$ws.Range("E119").Value = 50

# row 120: Kirguistan
$ws.Range("A120").Value = 'Kirguistan'
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 0
$ws.Range("E120").Value = 44
$ws.Range("F120").Value = 0
$ws.Range("H120").Value = 0

# row 121: Banglades
$ws.Range("A121").Value = 'Banglades'
$ws.Range("B121").Value = 44
$ws.Range("C121").Value = 5
$ws.Range("D121").Value = 11
$ws.Range("E121").Value = 28
$ws.Range("F121").Value = 1
$ws.Range("H121").Value = 5

# row 122: Bolivia
$ws.Range("A122").Value = 'Bolivia'
$ws.Range("B122").Value = 43
$ws.Range("C122").Value = 11
$ws.Range("E122").Value = 43

# row 132: Aruba
$ws.Range("A132").Value = 'Aruba'
$ws.Range("C132").Value = 9
$ws.Range("D132").Value = 1
$ws.Range("E132").Value = 27

# row 133: Guayana Francesa
$ws.Range("A133").Value = 'Guayana Francesa'
$ws.Range("B133").Value = 28
$ws.Range("D133").Value = 6
$ws.Range("E133").Value = 22
$ws.Range("H133").Value = 0

# row 134: Jamaica
$ws.Range("A134").Value = 'Jamaica'
$ws.Range("B134").Value = 26
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 2
$ws.Range("E134").Value = 23
$ws.Range("H134").Value = 1

# row 135: Isla de Man
$ws.Range("A135").Value = 'Isla de Man'
$ws.Range("B135").Value = 25
$ws.Range("C135").Value = 2
$ws.Range("D135").Value = 0
$ws.Range("E135").Value = 25
$ws.Range("H135").Value = 0

# row 136: Guatemala
$ws.Range("A136").Value = 'Guatemala'
$ws.Range("B136").Value = 25
$ws.Range("C136").Value = 1
$ws.Range("D136").Value = 4
$ws.Range("E136").Value = 20
$ws.Range("H136").Value = 1

# row 137: Madagascar
$ws.Range("A137").Value = 'Madagascar'
$ws.Range("C137").Value = 4
$ws.Range("D137").Value = 0
$ws.Range("E137").Value = 23

# row 138: Togo
$ws.Range("A138").Value = 'Togo'
$ws.Range("B138").Value = 23
$ws.Range("E138").Value = 22

# row 145: El Salvador
$ws.Range("A145").Value = 'El Salvador'
$ws.Range("C145").Value = 4

# row 146: Tanzania
$ws.Range("A146").Value = 'Tanzania'
$ws.Range("C146").Value = 0

# row 148: Etiopia
$ws.Range("A148").Value = 'Etiopia'
$ws.Range("C148").Value = 0

# row 149: Guinea Ecuatorial
$ws.Range("A149").Value = 'Guinea Ecuatorial'
$ws.Range("C149").Value = 3

# row 150: Mongolia
$ws.Range("A150").Value = 'Mongolia'
$ws.Range("C150").Value = 1

# row 151: San Martin (Parte Francesa)
$ws.Range("A151").Value = 'San Martin (Parte Francesa)'

# row 152: Republica de Yibuti
$ws.Range("A152").Value = 'Republica de Yibuti'
$ws.Range("C152").Value = 0

# row 153: Dominica
$ws.Range("A153").Value = 'Dominica'

# row 156: Surinam
$ws.Range("A156").Value = 'Surinam'

# row 157: Haiti
$ws.Range("A157").Value = 'Haiti'

# row 160: Antigua y Barbuda
$ws.Range("A160").Value = 'Antigua y Barbuda'
$ws.Range("C160").Value = 4

# row 161: Granada
$ws.Range("A161").Value = 'Granada'
$ws.Range("C161").Value = 6

# row 162: Mozambique
$ws.Range("A162").Value = 'Mozambique'
$ws.Range("C162").Value = 2

# row 163: Seychelles
$ws.Range("A163").Value = 'Seychelles'
$ws.Range("C163").Value = 0

# row 165: Eritrea
$ws.Range("A165").Value = 'Eritrea'
$ws.Range("C165").Value = 2

# row 166: Benin
$ws.Range("A166").Value = 'Benin'
$ws.Range("C166").Value = 0

# row 167: Suazilandia
$ws.Range("A167").Value = 'Suazilandia'

# row 168: Laos
$ws.Range("A168").Value = 'Laos'
$ws.Range("C168").Value = 3

# row 171: Birmania
$ws.Range("A171").Value = 'Birmania'
$ws.Range("C171").Value = 2

# row 172: Siria
$ws.Range("A172").Value = 'Siria'

# row 173: Montserrat
$ws.Range("A173").Value = 'Montserrat'
$ws.Range("C173").Value = 4
$ws.Range("E173").Value = 5
$ws.Range("H173").Value = 0

# row 174: Fiyi
$ws.Range("A174").Value = 'Fiyi'
$ws.Range("B174").Value = 5
$ws.Range("E174").Value = 5

# row 175: Guyana
$ws.Range("A175").Value = 'Guyana'
$ws.Range("B175").Value = 5
$ws.Range("C175").Value = 0
$ws.Range("H175").Value = 1

# row 176: Guinea
$ws.Range("A176").Value = 'Guinea'

# row 177: Mali
$ws.Range("A177").Value = 'Mali'
$ws.Range("C177").Value = 2

# row 178: Congo
$ws.Range("A178").Value = 'Congo'

# row 179: Angola
$ws.Range("A179").Value = 'Angola'
$ws.Range("C179").Value = 1
$ws.Range("E179").Value = 4
$ws.Range("H179").Value = 0

# row 180: Santa Sede
$ws.Range("A180").Value = 'Santa Sede'
$ws.Range("B180").Value = 4
$ws.Range("E180").Value = 4

# row 181: Cabo Verde
$ws.Range("A181").Value = 'Cabo Verde'
$ws.Range("B181").Value = 4
$ws.Range("C181").Value = 0
$ws.Range("H181").Value = 1

# row 182: Republica del Chad
$ws.Range("A182").Value = 'Republica del Chad'

# row 183: San Bartolome
$ws.Range("A183").Value = 'San Bartolome'

# row 184: Liberia
$ws.Range("A184").Value = 'Liberia'

# row 185: Santa Lucia
$ws.Range("A185").Value = 'Santa Lucia'

# row 186: Republica de Africa Central
$ws.Range("A186").Value = 'Republica de Africa Central'

# row 187: Mauritania
$ws.Range("A187").Value = 'Mauritania'
$ws.Range("C187").Value = 1

# row 188: San Martin (Parte Holandesa)
$ws.Range("A188").Value = 'San Martin (Parte Holandesa)'
$ws.Range("E188").Value = 3
$ws.Range("H188").Value = 0

# row 190: Nepal
$ws.Range("A190").Value = 'Nepal'
$ws.Range("D190").Value = 1
$ws.Range("H190").Value = 0

# row 191: Zimbabue
$ws.Range("A191").Value = 'Zimbabue'
$ws.Range("D191").Value = 0
$ws.Range("H191").Value = 1

# row 192: Gambia
$ws.Range("A192").Value = 'Gambia'
$ws.Range("B192").Value = 3
$ws.Range("H192").Value = 1

# row 193: Anguila
$ws.Range("A193").Value = 'Anguila'
$ws.Range("C193").Value = 2

# row 194: Nicaragua
$ws.Range("A194").Value = 'Nicaragua'
$ws.Range("C194").Value = 0

# row 195: Butan
$ws.Range("A195").Value = 'Butan'

# row 196: Islas Turcas y Caicos
$ws.Range("A196").Value = 'Islas Turcas y Caicos'
$ws.Range("C196").Value = 1

# row 197: Belice
$ws.Range("A197").Value = 'Belice'

# row 198: Guinea-Bisau
$ws.Range("A198").Value = 'Guinea-Bisau'

# row 199: San Cristobal y Nieves
$ws.Range("A199").Value = 'San Cristobal y Nieves'
$ws.Range("C199").Value = 0

# row 200: Islas Virgenes Britanicas
$ws.Range("A200").Value = 'Islas Virgenes Britanicas'
$ws.Range("C200").Value = 0

# row 201: Somalia
$ws.Range("A201").Value = 'Somalia'
$ws.Range("B201").Value = 2
$ws.Range("C201").Value = 1
$ws.Range("E201").Value = 2

# row 202: Libia
$ws.Range("A202").Value = 'Libia'

# row 203: San Vicente y las Granadinas
$ws.Range("A203").Value = 'San Vicente y las Granadinas'

# row 204: Timor Oriental
$ws.Range("A204").Value = 'Timor Oriental'

# row 205: Papua Nueva Guinea
$ws.Range("A205").Value = 'Papua Nueva Guinea'
$ws.Range("B205").Value = 1
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 0
$ws.Range("E205").Value = 1
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

